$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 117  # was 115
$ws1.Range("F4").Value = 9250  # was 9248
$ws1.Range("F6").Value = 66  # was 65
$ws1.Range("F7").Value = 1959  # was 1958
$ws1.Range("F8").Value = 6389  # was 6387
$ws1.Range("F11").Value = 9745  # was 9744
$ws1.Range("F12").Value = 11064  # was 11061
$ws1.Range("F15").Value = 4911  # was 4910
$ws1.Range("F18").Value = 0  # was 95
$ws1.Range("F19").Value = 330  # was 329
$ws1.Range("F20").Value = 176  # was 175
$ws1.Range("F21").Value = 1331  # was 1330
$ws1.Range("F25").Value = 854  # was 853
$ws1.Range("F33").Value = 1734  # was 1735
$ws1.Range("F36").Value = 46  # was 45
$ws1.Range("F37").Value = 913  # was 912
$ws1.Range("F40").Value = 3307  # was 3304
$ws1.Range("F47").Value = 236  # was 235
$ws1.Range("F49").Value = 4202  # was 4200

# Sheet "本地生活" (Local Life) - sheet3
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5889  # was 5888

# Sheet "全部类型" (All Types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 9250  # was 9248
$ws4.Range("F4").Value = 66  # was 65
$ws4.Range("F7").Value = 6389  # was 6387
$ws4.Range("F9").Value = 9745  # was 9744
$ws4.Range("F10").Value = 11064  # was 11061
$ws4.Range("F14").Value = 4911  # was 4910
$ws4.Range("F18").Value = 176  # was 175
$ws4.Range("F20").Value = 1331  # was 1330
$ws4.Range("F24").Value = 854  # was 853
$ws4.Range("F31").Value = 1734  # was 1735
$ws4.Range("F36").Value = 46  # was 45
$ws4.Range("F37").Value = 913  # was 912
$ws4.Range("F48").Value = 236  # was 235
